$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay stored as text (it mixes thousands-dot
# formatted numbers, scientific-looking decimals and subscript notation),
# matching how the source data already stores these values as strings.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.193.49"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.578.87"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.72"
$ws.Range("E5").Value = "  -3.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.76"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.110"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.36"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.041.07"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.010.04"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000156"
$ws.Range("E16").Value = "  +5.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.563.45"
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.32"
$ws.Range("E18").Value = "  +4.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.78"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.29"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.23"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.66"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "560.18"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0840"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.74"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.20"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "166.57"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.410"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.45"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "166.87"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.55"
$ws.Range("E42").Value = "  -1.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.94"
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.53"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0582"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.628"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.03"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0251"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.94"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0233"
$ws.Range("E51").Value = "  +17.71%  "
